$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Referencias" row (row 13) with three new shared strings.
$ws.Range("A13").Value = "Referencias"
$ws.Range("B13").Value = "pag 99"
$ws.Range("C13").Value = "datasheet atmega32"

# New font (red) applied via a new cell style (cellXf) to the new row.
$ws.Range("A13:C13").Font.Color = 255

# Page orientation set to portrait (adds pageSetup element).
$ws.PageSetup.Orientation = 1

# Update the selected cell, matching the author's last cursor position.
$ws.Range("C18").Select()
